$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record (row) was inserted into the weekly log as row 4,
# pushing every existing record (previously rows 4..87) down by one row
# (to rows 5..88). Insert a whole row at position 4 so everything below
# shifts down, then populate the newly-inserted row with the new record.
$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = 5
$ws.Range("B4").Value = "Macroferia Regional de Talca"
$ws.Range("C4").Value = "Maule"
$ws.Range("D4").Value = 44552
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 100112022
$ws.Range("G4").Value = "Arveja Verde"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 16000
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 16000
$ws.Range("N4").Value = "$/saco 25 kilos"
$ws.Range("O4").Value = "Carahue"
$ws.Range("P4").Value = 640
$ws.Range("Q4").Value = 25
$ws.Range("R4").Value = "Hortaliza"
